$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Column C now holds term-accession values and needs its own explicit width.
$ws.Columns("C").ColumnWidth = 10.92

# --- Fix up formatting before values change ---
# C13 becomes a new wrapped-text cell (same look as the old E13/E14 accession cells).
$ws.Range("E14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# E13 loses the wrap-text look it had as an accession-number cell and becomes a plain cell.
$ws.Range("C12").Copy()
$ws.Range("E13").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 12 "Tags": harmonize tags, dropping the near-duplicate "protocol" /
#     "phenotyping " tags and shifting the remaining ones left. ---
$ws.Range("C12").Value = "Plant"
$ws.Range("D12").Value = " metadata "
$ws.Range("E12").Value = "assay"
$ws.Range("F12").Value = "MIAPPE"
$ws.Range("G12").ClearContents()
$ws.Range("H12").ClearContents()

# --- Row 13 "Tags Term Accession Number": give the remaining tags proper
#     term accession numbers. ---
$ws.Range("C13").Value = "NCIT:C14258"
$ws.Range("E13").Value = "OBI:0000070"
$ws.Range("B13").Value = "DPBO:1000224"
$ws.Rows("13").RowHeight = 28.8

# --- Row 14 "Tags Term Source REF": no longer needed now that accession
#     numbers are filled in directly. ---
$ws.Range("E14").ClearContents()

$ws.Range("D17").Select()
